$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (mirrors the refreshed crypto price feed)
$cells = [ordered]@{
    'D2' = '67.864.48'
    'D3' = '3.584.81'
    'E3' = '  +0.54%  '
    'E4' = '  -0.05%  '
    'D5' = '202.18'
    'E5' = '  +8.70%  '
    'D6' = '561.57'
    'E6' = '  -3.73%  '
    'D7' = '3.580.57'
    'E7' = '  +0.52%  '
    'D8' = '0.618'
    'E8' = '  +0.74%  '
    'E9' = '  -0.14%  '
    'D10' = '0.667'
    'E10' = '  -0.48%  '
    'D11' = '60.09'
    'E11' = '  +13.01%  '
    'E12' = '  +2.58%  '
    'E13' = '  +7.67%  '
    'D14' = '9.94'
    'E14' = '  +1.14%  '
    'D15' = '4.155.66'
    'E15' = '  +0.57%  '
    'D16' = '3.585.18'
    'E16' = '  +0.52%  '
    'E17' = '  +0.57%  '
    'D18' = '18.84'
    'E18' = '  +2.72%  '
    'D19' = '67.622.92'
    'E19' = '  +1.97%  '
    'D20' = '12.30'
    'E20' = '  +0.63%  '
    'E21' = '  +1.36%  '
    'D22' = '400.16'
    'E22' = '  +1.23%  '
    'D23' = '12.79'
    'E23' = '  +13.11%  '
    'D24' = '4.13'
    'E24' = '  -4.87%  '
    'D25' = '84.84'
    'E25' = '  -1.46%  '
    'D26' = '2.89'
    'E26' = '  -0.45%  '
    'D27' = '12.52'
    'E27' = '  +0.22%  '
    'D28' = '3.91'
    'E28' = '  +10.37%  '
    'E29' = '  +1.25%  '
    'D30' = '8.21'
    'E30' = '  +15.50%  '
    'D31' = '9.29'
    'E31' = '  +3.79%  '
    'D32' = '31.41'
    'E32' = '  +0.88%  '
    'D33' = '664.53'
    'E33' = '  +6.73%  '
    'D34' = '12.13'
    'E34' = '  -0.40%  '
    'E35' = '  -0.04%  '
    'D36' = '63.33'
    'E36' = '  +0.07%  '
    'D37' = '41.85'
    'E37' = '  +0.94%  '
    'D38' = '0.420'
    'E38' = '  +5.38%  '
    'E39' = '  +0.10%  '
    'D40' = '3.294.12'
    'E40' = '  +9.55%  '
    'D41' = '0.0₃0757'
    'E41' = '  -0.65%  '
    'E42' = '  +11.79%  '
    'D43' = '0.135'
    'E43' = '  +2.45%  '
    'D44' = '2.75'
    'E44' = '  +8.95%  '
    'D45' = '0.997'
    'E45' = '  -0.18%  '
    'D46' = '2.95'
    'E46' = '  +27.06%  '
    'D47' = '0.0415'
    'E47' = '  +1.49%  '
    'D48' = '2.73'
    'E48' = '  +10.54%  '
    'B49' = 'ApeXProtocol'
    'C49' = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
    'D49' = '3.12'
    'E49' = '  -0.33%  '
    'D50' = '0.130'
    'E50' = '  -0.33%  '
    'B51' = 'THORChain'
    'C51' = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
    'D51' = '8.76'
    'E51' = '  +2.17%  '
}

foreach ($ref in $cells.Keys) {
    $range = $ws.Range($ref)
    # Force text interpretation (so numeric-looking strings like '202.18' are not
    # converted to numbers), then drop the formatting footprint so the cell is left
    # exactly as it was - a plain, unstyled inline string cell.
    $range.NumberFormat = '@'
    $range.Value = $cells[$ref]
    $range.ClearFormats()
}
